# Scheduled runner update: refresh computed market-price/profit figures
# across the per-job Leve sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 610.4
$ws.Range("J33").Value = 698.4545000000001
$ws.Range("L33").Value = 698.4545000000001
$ws.Range("N33").Value = -1156.4545
$ws.Range("H82").Value = 4160.5
$ws.Range("I82").Value = 321
$ws.Range("K82").Value = 963
$ws.Range("M82").Value = -557
$ws.Range("H85").Value = 4160.5
$ws.Range("I85").Value = 321
$ws.Range("K85").Value = 963
$ws.Range("M85").Value = 441
$ws.Range("H100").Value = 4342.857
$ws.Range("I100").Value = 3816.6667
$ws.Range("J100").Value = 4737.5
$ws.Range("K100").Value = 3816.6667
$ws.Range("L100").Value = 4737.5
$ws.Range("M100").Value = -3275.6667
$ws.Range("N100").Value = -5819.5
$ws.Range("H113").Value = 5188.5713
$ws.Range("I113").Value = 4944
$ws.Range("J113").Value = 5800
$ws.Range("K113").Value = 4944
$ws.Range("L113").Value = 5800
$ws.Range("M113").Value = -1690
$ws.Range("N113").Value = -12308
$ws.Range("H116").Value = 6614.5454
$ws.Range("I116").Value = 7568
$ws.Range("J116").Value = 5820
$ws.Range("K116").Value = 7568
$ws.Range("L116").Value = 5820
$ws.Range("M116").Value = -4126
$ws.Range("N116").Value = -12704

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1545.3334
$ws.Range("I45").Value = 1171.7715
$ws.Range("J45").Value = 2362.5
$ws.Range("K45").Value = 1171.7715
$ws.Range("L45").Value = 2362.5
$ws.Range("M45").Value = -794.7715000000001
$ws.Range("N45").Value = -3116.5
$ws.Range("H61").Value = 5973.304
$ws.Range("I61").Value = 5818.6
$ws.Range("K61").Value = 5818.6
$ws.Range("M61").Value = -5606.6
$ws.Range("H110").Value = 2291.7144
$ws.Range("I110").Value = 707.1
$ws.Range("J110").Value = 6253.25
$ws.Range("K110").Value = 707.1
$ws.Range("L110").Value = 6253.25
$ws.Range("M110").Value = 1337.9
$ws.Range("N110").Value = -10343.25
$ws.Range("H132").Value = 1621.7246
$ws.Range("I132").Value = 1236.0328
$ws.Range("J132").Value = 4562.625
$ws.Range("K132").Value = 3708.0984
$ws.Range("L132").Value = 13687.875
$ws.Range("M132").Value = -1178.0984
$ws.Range("N132").Value = -18747.875
$ws.Range("H136").Value = 5973.304
$ws.Range("I136").Value = 5818.6
$ws.Range("K136").Value = 17455.8
$ws.Range("M136").Value = -14905.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2754
$ws.Range("I134").Value = 2107.9412
$ws.Range("J134").Value = 5499.75
$ws.Range("K134").Value = 6323.823600000001
$ws.Range("L134").Value = 16499.25
$ws.Range("M134").Value = -3788.823600000001
$ws.Range("N134").Value = -21569.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 16060005
$ws.Range("J12").Value = 75006
$ws.Range("L12").Value = 75006
$ws.Range("N12").Value = -75346
$ws.Range("H58").Value = 17245520
$ws.Range("I58").Value = 2957.4211
$ws.Range("K58").Value = 2957.4211
$ws.Range("M58").Value = -2754.4211
$ws.Range("H122").Value = 2276.7576
$ws.Range("I122").Value = 1628.1364
$ws.Range("J122").Value = 3574
$ws.Range("K122").Value = 4884.4092
$ws.Range("L122").Value = 10722
$ws.Range("M122").Value = -2434.4092
$ws.Range("N122").Value = -15622
$ws.Range("H132").Value = 2217.775
$ws.Range("I132").Value = 1770.8148
$ws.Range("J132").Value = 3146.077
$ws.Range("K132").Value = 5312.4444
$ws.Range("L132").Value = 9438.231
$ws.Range("M132").Value = -2782.4444
$ws.Range("N132").Value = -14498.231
$ws.Range("H134").Value = 2741.1428
$ws.Range("I134").Value = 967
$ws.Range("J134").Value = 5106.6665
$ws.Range("K134").Value = 2901
$ws.Range("L134").Value = 15319.9995
$ws.Range("M134").Value = -366
$ws.Range("N134").Value = -20389.9995
$ws.Range("H136").Value = 17245520
$ws.Range("I136").Value = 2957.4211
$ws.Range("K136").Value = 8872.263300000001
$ws.Range("M136").Value = -6322.263300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 1812.5
$ws.Range("I15").Value = 125
$ws.Range("K15").Value = 375
$ws.Range("M15").Value = -235
$ws.Range("H47").Value = 1956
$ws.Range("I47").Value = 66.666664
$ws.Range("K47").Value = 199.999992
$ws.Range("M47").Value = 231.000008
$ws.Range("H131").Value = 1069.5178
$ws.Range("J131").Value = 1121.4314
$ws.Range("L131").Value = 3364.2942
$ws.Range("N131").Value = -13444.2942

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4966.8335
$ws.Range("I122").Value = 5579.6
$ws.Range("J122").Value = 4731.154
$ws.Range("K122").Value = 16738.8
$ws.Range("L122").Value = 14193.462
$ws.Range("M122").Value = -14288.8
$ws.Range("N122").Value = -19093.462
$ws.Range("H132").Value = 3359.825
$ws.Range("I132").Value = 3066.4814
$ws.Range("J132").Value = 3969.077
$ws.Range("K132").Value = 9199.4442
$ws.Range("L132").Value = 11907.231
$ws.Range("M132").Value = -6669.4442
$ws.Range("N132").Value = -16967.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 42574.855
$ws.Range("J20").Value = 42574.855
$ws.Range("L20").Value = 42574.855
$ws.Range("N20").Value = -43026.855
$ws.Range("H122").Value = 3909.0908
$ws.Range("I122").Value = 2833.3333
$ws.Range("J122").Value = 5200
$ws.Range("K122").Value = 8499.999899999999
$ws.Range("L122").Value = 15600
$ws.Range("M122").Value = -6049.999899999999
$ws.Range("N122").Value = -20500
$ws.Range("H136").Value = 2327744.5
$ws.Range("I136").Value = 2858206
$ws.Range("J136").Value = 6975
$ws.Range("K136").Value = 8574618
$ws.Range("L136").Value = 20925
$ws.Range("M136").Value = -8572068
$ws.Range("N136").Value = -26025

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1522.2
$ws.Range("I81").Value = 1219.2
$ws.Range("J81").Value = 1724.2
$ws.Range("K81").Value = 2438.4
$ws.Range("L81").Value = 3448.4
$ws.Range("M81").Value = -1377.4
$ws.Range("N81").Value = -5570.4
$ws.Range("H84").Value = 1522.2
$ws.Range("I84").Value = 1219.2
$ws.Range("J84").Value = 1724.2
$ws.Range("K84").Value = 12192
$ws.Range("L84").Value = 17242
$ws.Range("M84").Value = -6888
$ws.Range("N84").Value = -27850
$ws.Range("H132").Value = 130672.14
$ws.Range("I132").Value = 157000.38
$ws.Range("J132").Value = 10314.5
$ws.Range("K132").Value = 471001.14
$ws.Range("L132").Value = 30943.5
$ws.Range("M132").Value = -468471.14
$ws.Range("N132").Value = -36003.5
$ws.Range("H138").Value = 29833.334
$ws.Range("J138").Value = 29833.334
$ws.Range("L138").Value = 29833.334
$ws.Range("N138").Value = -40113.334
